# Insert a new data row before current row 94 (the "Terminal La Palmera de La Serena" /
# "Zanahoria" weekly price table). This shifts all existing rows 94-207 down to 95-208,
# and the new row 94 is populated with a new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94, pushing rows 94..207 down to 95..208.
$ws.Rows("94").Insert()

# Populate the newly inserted row 94 with the new record's values.
$ws.Cells.Item(94, 1).Value = 8
$ws.Cells.Item(94, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(94, 3).Value = "Coquimbo"
$ws.Cells.Item(94, 4).Value = 44483
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = 100114013
$ws.Cells.Item(94, 7).Value = "Zanahoria"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 600
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 7000
$ws.Cells.Item(94, 13).Value = 6500
$ws.Cells.Item(94, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 325
$ws.Cells.Item(94, 17).Value = 20
$ws.Cells.Item(94, 18).Value = "Hortaliza"
